{"js": "// Append a line break followed by a new sentence to the very end of the\n// document body (the last paragraph, which holds the contact info /\n// \"pdf version \u2022 txt version \u2022 doc version \u2022 html version\" links).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a manual line break (\"textWrapping\") at the end of the paragraph.\nlastParagraph.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\n\n// Insert the new sentence after that break, still at the end of the\n// paragraph, as its own run.\nlastParagraph.insertText(\n  \"I prefer the email over other means of communication.\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Append a manual line break followed by a new sentence to the very end of\n# the document body (the last paragraph, which holds the contact info /\n# \"pdf version \u2022 txt version \u2022 doc version \u2022 html version\" links).\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Collapse(0)            # wdCollapseEnd\n$r.InsertBreak(6)         # wdLineBreak\n$r.Collapse(0)            # wdCollapseEnd\n$r.InsertAfter(\"I prefer the email over other means of communication.\")\n"}
